$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation for cells that would otherwise be
# auto-converted by Excel (date-looking string, zero-padded numeric id).
$ws.Range("B35").NumberFormat = "@"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("I35").NumberFormat = "@"

$ws.Range("A35").Value = 1583971200
$ws.Range("B35").Value = "2020-03-12"
$ws.Range("C35").Value = "03033"
$ws.Range("D35").Value = "SGBHD"
$ws.Range("E35").Value = 0.3
$ws.Range("F35").Value = 0.3
$ws.Range("G35").Value = 0.3
$ws.Range("H35").Value = 0.3
$ws.Range("I35").Value = "-"

# Reset to the default "Normal" style so the new row matches the
# unstyled look of the other data rows (keeps the text cell type).
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Style = "Normal"
$ws.Range("I35").Style = "Normal"
